# Scheduled-runner refresh of cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# "Durandal_Profits" leve tables. Values below are the updated snapshot
# pulled by the runner; row/column layout is unchanged.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1480.3
$ws.Range("I40").Value = 1360.6
$ws.Range("K40").Value = 1360.6
$ws.Range("M40").Value = -1185.6
$ws.Range("H62").Value = 2555.5
$ws.Range("I62").Value = 2456.5
$ws.Range("J62").Value = 3347.5
$ws.Range("K62").Value = 2456.5
$ws.Range("L62").Value = 3347.5
$ws.Range("M62").Value = -1832.5
$ws.Range("N62").Value = -4595.5
$ws.Range("H64").Value = 3235
$ws.Range("I64").Value = 2950
$ws.Range("J64").Value = 3520
$ws.Range("K64").Value = 2950
$ws.Range("L64").Value = 3520
$ws.Range("M64").Value = -2702
$ws.Range("N64").Value = -4016
$ws.Range("H65").Value = 2555.5
$ws.Range("I65").Value = 2456.5
$ws.Range("J65").Value = 3347.5
$ws.Range("K65").Value = 12282.5
$ws.Range("L65").Value = 16737.5
$ws.Range("M65").Value = -9162.5
$ws.Range("N65").Value = -22977.5
$ws.Range("H67").Value = 3235
$ws.Range("I67").Value = 2950
$ws.Range("J67").Value = 3520
$ws.Range("K67").Value = 2950
$ws.Range("L67").Value = 3520
$ws.Range("M67").Value = -2092
$ws.Range("N67").Value = -5236
$ws.Range("H76").Value = 2318121
$ws.Range("I76").Value = 2648688.2
$ws.Range("J76").Value = 4150
$ws.Range("K76").Value = 2648688.2
$ws.Range("L76").Value = 4150
$ws.Range("M76").Value = -2648373.2
$ws.Range("N76").Value = -4780
$ws.Range("H79").Value = 2318121
$ws.Range("I79").Value = 2648688.2
$ws.Range("J79").Value = 4150
$ws.Range("K79").Value = 2648688.2
$ws.Range("L79").Value = 4150
$ws.Range("M79").Value = -2647596.2
$ws.Range("N79").Value = -6334
$ws.Range("H80").Value = 2205.8918
$ws.Range("I80").Value = 675
$ws.Range("J80").Value = 4006.9412
$ws.Range("K80").Value = 2025
$ws.Range("L80").Value = 12020.8236
$ws.Range("M80").Value = -1027
$ws.Range("N80").Value = -14016.8236
$ws.Range("H83").Value = 2205.8918
$ws.Range("I83").Value = 675
$ws.Range("J83").Value = 4006.9412
$ws.Range("K83").Value = 6075
$ws.Range("L83").Value = 36062.4708
$ws.Range("M83").Value = -1083
$ws.Range("N83").Value = -46046.4708
$ws.Range("H100").Value = 12347206
$ws.Range("I100").Value = 15873865
$ws.Range("J100").Value = 3900
$ws.Range("K100").Value = 15873865
$ws.Range("L100").Value = 3900
$ws.Range("M100").Value = -15873324
$ws.Range("N100").Value = -4982
$ws.Range("H106").Value = 2246.5386
$ws.Range("I106").Value = 2183.75
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 2183.75
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -1552.75
$ws.Range("N106").Value = -4262
$ws.Range("H116").Value = 6057.846
$ws.Range("I116").Value = 7237.0527
$ws.Range("J116").Value = 2857.1428
$ws.Range("K116").Value = 7237.0527
$ws.Range("L116").Value = 2857.1428
$ws.Range("M116").Value = -3795.0527
$ws.Range("N116").Value = -9741.1428
$ws.Range("H129").Value = 974.3674
$ws.Range("I129").Value = 313.0909
$ws.Range("J129").Value = 1165.7894
$ws.Range("K129").Value = 939.2727
$ws.Range("L129").Value = 3497.3682
$ws.Range("M129").Value = 4060.7273
$ws.Range("N129").Value = -13497.3682
$ws.Range("H137").Value = 1901.7273
$ws.Range("I137").Value = 2034.6471
$ws.Range("J137").Value = 1449.8
$ws.Range("K137").Value = 6103.9413
$ws.Range("L137").Value = 4349.4
$ws.Range("M137").Value = -3553.9413
$ws.Range("N137").Value = -9449.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 50051576
$ws.Range("I132").Value = 66667704
$ws.Range("J132").Value = 203202.4
$ws.Range("K132").Value = 200003112
$ws.Range("L132").Value = 609607.2
$ws.Range("M132").Value = -200000582
$ws.Range("N132").Value = -614667.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1151.8518
$ws.Range("I99").Value = 853.1579
$ws.Range("J99").Value = 1861.25
$ws.Range("K99").Value = 853.1579
$ws.Range("L99").Value = 1861.25
$ws.Range("M99").Value = 644.8421
$ws.Range("N99").Value = -4857.25
$ws.Range("H134").Value = 11378.053
$ws.Range("I134").Value = 3898.125
$ws.Range("J134").Value = 51271
$ws.Range("K134").Value = 11694.375
$ws.Range("L134").Value = 153813
$ws.Range("M134").Value = -9159.375
$ws.Range("N134").Value = -158883

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2146.75
$ws.Range("I31").Value = 1616.75
$ws.Range("J31").Value = 4001.75
$ws.Range("K31").Value = 1616.75
$ws.Range("L31").Value = 4001.75
$ws.Range("M31").Value = -1321.75
$ws.Range("N31").Value = -4591.75
$ws.Range("H34").Value = 2146.75
$ws.Range("I34").Value = 1616.75
$ws.Range("J34").Value = 4001.75
$ws.Range("K34").Value = 1616.75
$ws.Range("L34").Value = 4001.75
$ws.Range("M34").Value = -1414.75
$ws.Range("N34").Value = -4405.75
$ws.Range("H58").Value = 780.41174
$ws.Range("I58").Value = 678.7406999999999
$ws.Range("J58").Value = 1172.5714
$ws.Range("K58").Value = 678.7406999999999
$ws.Range("L58").Value = 1172.5714
$ws.Range("M58").Value = -475.7406999999999
$ws.Range("N58").Value = -1578.5714
$ws.Range("H136").Value = 780.41174
$ws.Range("I136").Value = 678.7406999999999
$ws.Range("J136").Value = 1172.5714
$ws.Range("K136").Value = 2036.2221
$ws.Range("L136").Value = 3517.7142
$ws.Range("M136").Value = 513.7779
$ws.Range("N136").Value = -8617.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 20000
$ws.Range("J139").Value = 20000
$ws.Range("L139").Value = 20000
$ws.Range("N139").Value = -30280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2618.7896
$ws.Range("I7").Value = 1839.7858
$ws.Range("J7").Value = 4800
$ws.Range("K7").Value = 1839.7858
$ws.Range("L7").Value = 4800
$ws.Range("M7").Value = -1727.7858
$ws.Range("N7").Value = -5024
$ws.Range("H94").Value = 39999
$ws.Range("J94").Value = 39999
$ws.Range("L94").Value = 39999
$ws.Range("N94").Value = -41351
$ws.Range("H126").Value = 2618.7896
$ws.Range("I126").Value = 1839.7858
$ws.Range("J126").Value = 4800
$ws.Range("K126").Value = 5519.357400000001
$ws.Range("L126").Value = 14400
$ws.Range("M126").Value = -3049.357400000001
$ws.Range("N126").Value = -19340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 79533.336
$ws.Range("J138").Value = 79533.336
$ws.Range("L138").Value = 79533.336
$ws.Range("N138").Value = -89813.336
